# Apply updated crypto price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.705.92"
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("D3").Value = "2.252.00"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'306.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "'95.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.519"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("D10").Value = "'34.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.48%  "
$ws.Range("D11").Value = "'0.0802"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("D14").Value = "2.392.96"
$ws.Range("E14").Value = "  +2.59%  "
$ws.Range("D15").Value = "2.594.68"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("D17").Value = "'13.61"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").Value = "44.477.63"
$ws.Range("D19").Value = "0.0₃0940"
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("E20").Value = "  -2.55%  "
$ws.Range("E21").Value = "  -3.10%  "
$ws.Range("D22").Value = "'65.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("D23").Value = "'238.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "'2.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "'1.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.13%  "
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("E27").Value = "  +3.96%  "
$ws.Range("E28").Value = "  -1.47%  "
$ws.Range("D29").Value = "'37.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.41%  "
$ws.Range("E30").Value = "  +1.11%  "
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D32").Value = "'148.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.09%  "
$ws.Range("D33").Value = "'0.0786"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.16%  "
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("E35").Value = "  +0.95%  "
$ws.Range("D36").Value = "'0.108"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'1.87"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.32%  "
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").Value = "'0.118"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.41%  "
$ws.Range("D39").Value = "'15.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.90%  "
$ws.Range("E40").Value = "  -5.58%  "
$ws.Range("D41").Value = "'3.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.67%  "
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").Value = "1.802.56"
$ws.Range("E44").Value = "  +3.27%  "
$ws.Range("E45").Value = "  +12.10%  "
$ws.Range("D46").Value = "'82.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("E47").Value = "  -1.74%  "
$ws.Range("D48").Value = "'98.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("D49").Value = "'4.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.20%  "
$ws.Range("D50").Value = "'69.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.84%  "
$ws.Range("D51").Value = "'54.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.95%  "
